$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Газпром")
Write-Host $ws.Name
